$d = $word.ActiveDocument
$d.Content.Find.Execute("94-59=", $true, $false, $false, $false, $false, $true, 1, $false, "1+31=", 2) | Out-Null
$d.Content.Find.Execute("6+20=", $true, $false, $false, $false, $false, $true, 1, $false, "81-45=", 2) | Out-Null
$d.Content.Find.Execute("83-79=", $true, $false, $false, $false, $false, $true, 1, $false, "35+58=", 2) | Out-Null
$d.Content.Find.Execute("60-13=", $true, $false, $false, $false, $false, $true, 1, $false, "84-38=", 2) | Out-Null
$d.Content.Find.Execute("14+61=", $true, $false, $false, $false, $false, $true, 1, $false, "28+18=", 2) | Out-Null
$d.Content.Find.Execute("32+43=", $true, $false, $false, $false, $false, $true, 1, $false, "83-20=", 2) | Out-Null
$d.Content.Find.Execute("41+32=", $true, $false, $false, $false, $false, $true, 1, $false, "62+35=", 2) | Out-Null
$d.Content.Find.Execute("31+38=", $true, $false, $false, $false, $false, $true, 1, $false, "11+46=", 2) | Out-Null
$d.Content.Find.Execute("53+29=", $true, $false, $false, $false, $false, $true, 1, $false, "50-44=", 2) | Out-Null
$d.Content.Find.Execute("26+57=", $true, $false, $false, $false, $false, $true, 1, $false, "50-46=", 2) | Out-Null
$d.Content.Find.Execute("4+79=", $true, $false, $false, $false, $false, $true, 1, $false, "20+31=", 2) | Out-Null
$d.Content.Find.Execute("96-61=", $true, $false, $false, $false, $false, $true, 1, $false, "14+4=", 2) | Out-Null
$d.Content.Find.Execute("44-16=", $true, $false, $false, $false, $false, $true, 1, $false, "3+71=", 2) | Out-Null
$d.Content.Find.Execute("26+54=", $true, $false, $false, $false, $false, $true, 1, $false, "40+0=", 2) | Out-Null
$d.Content.Find.Execute("36+51=", $true, $false, $false, $false, $false, $true, 1, $false, "7+91=", 2) | Out-Null
$d.Content.Find.Execute("10+78=", $true, $false, $false, $false, $false, $true, 1, $false, "81-43=", 2) | Out-Null
$d.Content.Find.Execute("62+29=", $true, $false, $false, $false, $false, $true, 1, $false, "71-47=", 2) | Out-Null
$d.Content.Find.Execute("36+46=", $true, $false, $false, $false, $false, $true, 1, $false, "71+0=", 2) | Out-Null
$d.Content.Find.Execute("81+3=", $true, $false, $false, $false, $false, $true, 1, $false, "36-8=", 2) | Out-Null
$d.Content.Find.Execute("97-91=", $true, $false, $false, $false, $false, $true, 1, $false, "75-57=", 2) | Out-Null
$d.Content.Find.Execute("24+61=", $true, $false, $false, $false, $false, $true, 1, $false, "66-2=", 2) | Out-Null
$d.Content.Find.Execute("11+24=", $true, $false, $false, $false, $false, $true, 1, $false, "78-1=", 2) | Out-Null
$d.Content.Find.Execute("41+48=", $true, $false, $false, $false, $false, $true, 1, $false, "35+43=", 2) | Out-Null
$d.Content.Find.Execute("12+12=", $true, $false, $false, $false, $false, $true, 1, $false, "5-4=", 2) | Out-Null
$d.Content.Find.Execute("28-13=", $true, $false, $false, $false, $false, $true, 1, $false, "28-9=", 2) | Out-Null
$d.Content.Find.Execute("57-6=", $true, $false, $false, $false, $false, $true, 1, $false, "59-57=", 2) | Out-Null
$d.Content.Find.Execute("12+48=", $true, $false, $false, $false, $false, $true, 1, $false, "63+27=", 2) | Out-Null
$d.Content.Find.Execute("11+33=", $true, $false, $false, $false, $false, $true, 1, $false, "66-17=", 2) | Out-Null
$d.Content.Find.Execute("23+24=", $true, $false, $false, $false, $false, $true, 1, $false, "87-40=", 2) | Out-Null
$d.Content.Find.Execute("73-70=", $true, $false, $false, $false, $false, $true, 1, $false, "22+6=", 2) | Out-Null
$d.Content.Find.Execute("53+46=", $true, $false, $false, $false, $false, $true, 1, $false, "76-17=", 2) | Out-Null
$d.Content.Find.Execute("97-75=", $true, $false, $false, $false, $false, $true, 1, $false, "21-3=", 2) | Out-Null
$d.Content.Find.Execute("82-20=", $true, $false, $false, $false, $false, $true, 1, $false, "59+6=", 2) | Out-Null
$d.Content.Find.Execute("91-51=", $true, $false, $false, $false, $false, $true, 1, $false, "97-34=", 2) | Out-Null
$d.Content.Find.Execute("48+41=", $true, $false, $false, $false, $false, $true, 1, $false, "94+4=", 2) | Out-Null
$d.Content.Find.Execute("38-27=", $true, $false, $false, $false, $false, $true, 1, $false, "4+65=", 2) | Out-Null
$d.Content.Find.Execute("23-14=", $true, $false, $false, $false, $false, $true, 1, $false, "65-18=", 2) | Out-Null
$d.Content.Find.Execute("84-46=", $true, $false, $false, $false, $false, $true, 1, $false, "56-21=", 2) | Out-Null
$d.Content.Find.Execute("17+38=", $true, $false, $false, $false, $false, $true, 1, $false, "35+5=", 2) | Out-Null
$d.Content.Find.Execute("97-6=", $true, $false, $false, $false, $false, $true, 1, $false, "45+25=", 2) | Out-Null
$d.Content.Find.Execute("77-25=", $true, $false, $false, $false, $false, $true, 1, $false, "38+26=", 2) | Out-Null
$d.Content.Find.Execute("74-72=", $true, $false, $false, $false, $false, $true, 1, $false, "59-18=", 2) | Out-Null
$d.Content.Find.Execute("42-17=", $true, $false, $false, $false, $false, $true, 1, $false, "82+10=", 2) | Out-Null
$d.Content.Find.Execute("60-31=", $true, $false, $false, $false, $false, $true, 1, $false, "50-33=", 2) | Out-Null
$d.Content.Find.Execute("55-26=", $true, $false, $false, $false, $false, $true, 1, $false, "13+54=", 2) | Out-Null
$d.Content.Find.Execute("29-3=", $true, $false, $false, $false, $false, $true, 1, $false, "22+55=", 2) | Out-Null
$d.Content.Find.Execute("9+76=", $true, $false, $false, $false, $false, $true, 1, $false, "26+29=", 2) | Out-Null
$d.Content.Find.Execute("58+39=", $true, $false, $false, $false, $false, $true, 1, $false, "55-15=", 2) | Out-Null
$d.Content.Find.Execute("6+43=", $true, $false, $false, $false, $false, $true, 1, $false, "12-1=", 2) | Out-Null
$d.Content.Find.Execute("26+39=", $true, $false, $false, $false, $false, $true, 1, $false, "69+28=", 2) | Out-Null
$d.Content.Find.Execute("26+51=", $true, $false, $false, $false, $false, $true, 1, $false, "49-31=", 2) | Out-Null
$d.Content.Find.Execute("93-35=", $true, $false, $false, $false, $false, $true, 1, $false, "62-51=", 2) | Out-Null
$d.Content.Find.Execute("96-57=", $true, $false, $false, $false, $false, $true, 1, $false, "33-32=", 2) | Out-Null
$d.Content.Find.Execute("85-5=", $true, $false, $false, $false, $false, $true, 1, $false, "42+31=", 2) | Out-Null
$d.Content.Find.Execute("23+39=", $true, $false, $false, $false, $false, $true, 1, $false, "85-28=", 2) | Out-Null
$d.Content.Find.Execute("19+15=", $true, $false, $false, $false, $false, $true, 1, $false, "6+49=", 2) | Out-Null
$d.Content.Find.Execute("62-39=", $true, $false, $false, $false, $false, $true, 1, $false, "92-53=", 2) | Out-Null
$d.Content.Find.Execute("95-65=", $true, $false, $false, $false, $false, $true, 1, $false, "85-62=", 2) | Out-Null
$d.Content.Find.Execute("44-9=", $true, $false, $false, $false, $false, $true, 1, $false, "90-77=", 2) | Out-Null
$d.Content.Find.Execute("45-14=", $true, $false, $false, $false, $false, $true, 1, $false, "68-63=", 2) | Out-Null
$d.Content.Find.Execute("76-56=", $true, $false, $false, $false, $false, $true, 1, $false, "77-20=", 2) | Out-Null
$d.Content.Find.Execute("47-7=", $true, $false, $false, $false, $false, $true, 1, $false, "45+13=", 2) | Out-Null
$d.Content.Find.Execute("61-51=", $true, $false, $false, $false, $false, $true, 1, $false, "39+11=", 2) | Out-Null
$d.Content.Find.Execute("8+32=", $true, $false, $false, $false, $false, $true, 1, $false, "74-46=", 2) | Out-Null
$d.Content.Find.Execute("46-39=", $true, $false, $false, $false, $false, $true, 1, $false, "5+73=", 2) | Out-Null
$d.Content.Find.Execute("31+66=", $true, $false, $false, $false, $false, $true, 1, $false, "65+21=", 2) | Out-Null
$d.Content.Find.Execute("66-46=", $true, $false, $false, $false, $false, $true, 1, $false, "28-6=", 2) | Out-Null
$d.Content.Find.Execute("80+7=", $true, $false, $false, $false, $false, $true, 1, $false, "49-21=", 2) | Out-Null
$d.Content.Find.Execute("14+14=", $true, $false, $false, $false, $false, $true, 1, $false, "1+16=", 2) | Out-Null
$d.Content.Find.Execute("99-55=", $true, $false, $false, $false, $false, $true, 1, $false, "27+11=", 2) | Out-Null
$d.Content.Find.Execute("42-13=", $true, $false, $false, $false, $false, $true, 1, $false, "27+15=", 2) | Out-Null
$d.Content.Find.Execute("90-16=", $true, $false, $false, $false, $false, $true, 1, $false, "17+60=", 2) | Out-Null
$d.Content.Find.Execute("35-25=", $true, $false, $false, $false, $false, $true, 1, $false, "48-44=", 2) | Out-Null
$d.Content.Find.Execute("49+16=", $true, $false, $false, $false, $false, $true, 1, $false, "57-30=", 2) | Out-Null
$d.Content.Find.Execute("97-81=", $true, $false, $false, $false, $false, $true, 1, $false, "65-24=", 2) | Out-Null
$d.Content.Find.Execute("52+41=", $true, $false, $false, $false, $false, $true, 1, $false, "47-24=", 2) | Out-Null
$d.Content.Find.Execute("12+31=", $true, $false, $false, $false, $false, $true, 1, $false, "0+90=", 2) | Out-Null
$d.Content.Find.Execute("52-33=", $true, $false, $false, $false, $false, $true, 1, $false, "32+41=", 2) | Out-Null
$d.Content.Find.Execute("66-36=", $true, $false, $false, $false, $false, $true, 1, $false, "35+3=", 2) | Out-Null
$d.Content.Find.Execute("38+38=", $true, $false, $false, $false, $false, $true, 1, $false, "81-5=", 2) | Out-Null
$d.Content.Find.Execute("64-52=", $true, $false, $false, $false, $false, $true, 1, $false, "29+50=", 2) | Out-Null
$d.Content.Find.Execute("51-32=", $true, $false, $false, $false, $false, $true, 1, $false, "3+37=", 2) | Out-Null
$d.Content.Find.Execute("85-24=", $true, $false, $false, $false, $false, $true, 1, $false, "87-41=", 2) | Out-Null
$d.Content.Find.Execute("99-67=", $true, $false, $false, $false, $false, $true, 1, $false, "22+11=", 2) | Out-Null
$d.Content.Find.Execute("3+64=", $true, $false, $false, $false, $false, $true, 1, $false, "62-27=", 2) | Out-Null
$d.Content.Find.Execute("89-32=", $true, $false, $false, $false, $false, $true, 1, $false, "69-54=", 2) | Out-Null
$d.Content.Find.Execute("22-17=", $true, $false, $false, $false, $false, $true, 1, $false, "60+20=", 2) | Out-Null
$d.Content.Find.Execute("9+74=", $true, $false, $false, $false, $false, $true, 1, $false, "73-32=", 2) | Out-Null
$d.Content.Find.Execute("60-58=", $true, $false, $false, $false, $false, $true, 1, $false, "75+0=", 2) | Out-Null
$d.Content.Find.Execute("97-65=", $true, $false, $false, $false, $false, $true, 1, $false, "56+21=", 2) | Out-Null
$d.Content.Find.Execute("6+90=", $true, $false, $false, $false, $false, $true, 1, $false, "81-45=", 2) | Out-Null
$d.Content.Find.Execute("34-9=", $true, $false, $false, $false, $false, $true, 1, $false, "32+63=", 2) | Out-Null
$d.Content.Find.Execute("28-27=", $true, $false, $false, $false, $false, $true, 1, $false, "18+55=", 2) | Out-Null
$d.Content.Find.Execute("80+1=", $true, $false, $false, $false, $false, $true, 1, $false, "88-78=", 2) | Out-Null
$d.Content.Find.Execute("86-22=", $true, $false, $false, $false, $false, $true, 1, $false, "94-15=", 2) | Out-Null
$d.Content.Find.Execute("99-61=", $true, $false, $false, $false, $false, $true, 1, $false, "78-26=", 2) | Out-Null
$d.Content.Find.Execute("69-49=", $true, $false, $false, $false, $false, $true, 1, $false, "18+14=", 2) | Out-Null
$d.Content.Find.Execute("44-29=", $true, $false, $false, $false, $false, $true, 1, $false, "1+70=", 2) | Out-Null
$d.Content.Find.Execute("42-6=", $true, $false, $false, $false, $false, $true, 1, $false, "12+55=", 2) | Out-Null
$d.Content.Find.Execute("90-39=", $true, $false, $false, $false, $false, $true, 1, $false, "8+31=", 2) | Out-Null
